$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Degree name එක හරියට Mention කරලා තියෙනවද?"
$ws.Range("A8").Value = "Specialization area එක Mention කරලා ද?"
$ws.Range("A11").Value = "Add current GPA, if you have good results"
$ws.Range("A6").Value = "(Bachelor of Information and Communication Technology (Hons))"
$ws.Range("A9").Value = "(Software Technology/Network Technology/Multimedia Technology)"
$ws.Range("A13").Value = "CV එකේ Pages එකකට වඩා තිබීම"

$ws.Range("A13").Select()
